$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 7).Value = 'eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Cells.Item(4, 7).Value = 'eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Cells.Item(5, 7).Value = 'eman.tantawi@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, hananragab@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
$ws.Cells.Item(6, 7).Value = 'eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg'
$ws.Cells.Item(7, 7).Value = 'eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Cells.Item(8, 7).Value = 'eman.tantawi@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, backup@backdoor.com, asmaa.reda@med.asu.edu.eg, manar.montaser@med.asu.edu.eg'
$ws.Cells.Item(9, 7).Value = 'mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Cells.Item(10, 7).Value = 'alshimaa.atef@med.asu.edu.egm, shaimaa.ahmed@med.asu.edu.eg, drsarahmw@gmail.com, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Cells.Item(13, 7).Value = 'mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg'
$ws.Cells.Item(16, 7).Value = 'norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Cells.Item(17, 7).Value = 'marwa_mustafa@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg'
$ws.Cells.Item(19, 7).Value = 'Sarah.Mahdy@med.asu.edu.eg, mariam.e.mohammed@med.asu.edu.eg'
$ws.Cells.Item(24, 7).Value = 'Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Cells.Item(25, 7).Value = 'Remon.Matta@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Cells.Item(26, 7).Value = 'youstina.magdy@med.asu.edu.eg, gigi@med.asu.edu.eg'
$ws.Cells.Item(27, 7).Value = 'Eman_mohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(28, 7).Value = 'nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(29, 7).Value = 'Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Cells.Item(30, 7).Value = 'eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Cells.Item(31, 7).Value = 'eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Cells.Item(32, 7).Value = 'eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Cells.Item(33, 7).Value = 'eman.tantawi@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, hananragab@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
$ws.Cells.Item(34, 7).Value = 'eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg'
$ws.Cells.Item(35, 7).Value = 'eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Cells.Item(36, 7).Value = 'eman.tantawi@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, backup@backdoor.com, asmaa.reda@med.asu.edu.eg, manar.montaser@med.asu.edu.eg'
$ws.Cells.Item(37, 7).Value = 'mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Cells.Item(38, 7).Value = 'alshimaa.atef@med.asu.edu.egm, shaimaa.ahmed@med.asu.edu.eg, drsarahmw@gmail.com, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Cells.Item(41, 7).Value = 'mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg'
$ws.Cells.Item(44, 7).Value = 'norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Cells.Item(45, 7).Value = 'marwa_mustafa@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg'
$ws.Cells.Item(47, 7).Value = 'Sarah.Mahdy@med.asu.edu.eg, mariam.e.mohammed@med.asu.edu.eg'
$ws.Cells.Item(52, 7).Value = 'Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Cells.Item(53, 7).Value = 'Remon.Matta@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Cells.Item(54, 7).Value = 'youstina.magdy@med.asu.edu.eg, gigi@med.asu.edu.eg'
$ws.Cells.Item(55, 7).Value = 'Eman_mohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(56, 7).Value = 'nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(57, 7).Value = 'Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
